$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.372.39"
$ws.Range("E2").Value = "  +3.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.820.16"
$ws.Range("E3").Value = "  +4.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.22"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4455"
$ws.Range("E7").Value = "  +5.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3725"
$ws.Range("E8").Value = "  +3.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.97"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07764"
$ws.Range("E10").Value = "  +4.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.137"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("E12").Value = "  +3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.327"
$ws.Range("E14").Value = "  +4.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.626"
$ws.Range("E15").Value = "  +6.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.845.02"
$ws.Range("E16").Value = "  +6.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.58"
$ws.Range("E17").Value = "  +6.96%  "
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("E19").Value = "  +8.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.59"
$ws.Range("E21").Value = "  +4.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.261"
$ws.Range("E22").Value = "  +2.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.414.76"
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.224"
$ws.Range("E25").Value = "  -6.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.87"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.79"
$ws.Range("E27").Value = "  +5.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.051.11"
$ws.Range("E28").Value = "  +6.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.328"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.08"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.213"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.965"
$ws.Range("E32").Value = "  +5.33%  "
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.652"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.23"
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2199"
$ws.Range("E37").Value = "  +2.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.222"
$ws.Range("E38").Value = "  +3.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6639"
$ws.Range("E39").Value = "  +4.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06250"
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.202"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.192"
$ws.Range("E42").Value = "  +2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.426"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").Value = "  +3.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6172"
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.779"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.050"
$ws.Range("E48").Value = "  +4.86%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.30"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E50").Value = "  +6.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07012"
$ws.Range("E51").Value = "  +2.35%  "
